$wb = $excel.ActiveWorkbook

$sheetNames = @("10_trees", "500_trees", "1600_trees")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Add new "average" row at row 105
    $ws.Range("D105").Value2 = "average"

    if ($name -eq "10_trees") {
        $ws.Range("E105").Formula = "=AVERAGE(E2:E103,'500_trees'!E2:E103,'1600_trees'!E2:E103)"
    } else {
        $ws.Range("E105").Formula = "='10_trees'!E105"
    }

    # Update every existing K-column formula (SUM(E#,G#,H#)) to reference the
    # new average cell ($E$105) instead of the row's own E value.
    for ($r = 2; $r -le 103; $r++) {
        $kCell = $ws.Range("K$r")
        if ($kCell.HasFormula) {
            $kCell.Formula = "=SUM(`$E`$105,G$r,H$r)"
        }
    }
}

$excel.CalculateFull()
